$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9613333333333334
$ws.Range("C2").Value = 0.8246666666666667

$ws.Range("B3").Value = 0.958
$ws.Range("C3").Value = 0.8293333333333334

$ws.Range("B4").Value = 0.9606666666666667
$ws.Range("C4").Value = 0.7626666666666667

$ws.Range("B5").Value = 0.962
$ws.Range("C5").Value = 0.8386666666666667

$ws.Range("B6").Value = 0.9646666666666667
$ws.Range("C6").Value = 0.7353333333333333
